$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New header cells for columns I (I0) and J (IF)
$ws.Range("I1").Value = "I0"
$ws.Range("J1").Value = "IF"

# Match the header formatting used by the existing header row (e.g. H1: bold, bordered, centered)
$ws.Range("H1").Copy()
$ws.Range("I1:J1").PasteSpecial(-4122)

$I_VALUES = @(8,1,4,8,3,7,1,9,1,8,7,1,5,8,8,6,9,8,4,7,6,8,4,8,6,5,5,4,8,9,8,8,7,6,5,5,8,4,7,6,5,6,8,8,8,7,6,7,7,6,8,7,6,5,7,7,4,7,8,5,8,7,6,7,7,7,6,6)
$J_VALUES = @(9,1,6,8,4,7,1,9,1,8,7,2,6,8,8,6,9,8,5,7,6,8,5,8,7,6,6,4,8,9,8,8,7,7,5,6,8,5,7,6,6,6,8,8,8,7,7,7,7,6,8,8,7,6,7,7,5,8,8,5,8,7,6,8,7,8,6,6)

for ($r = 2; $r -le 69; $r++) {
    $idx = $r - 2
    $ws.Cells.Item($r, 9).Value = $I_VALUES[$idx]
    $ws.Cells.Item($r, 10).Value = $J_VALUES[$idx]
}

$wb.Save()
